# Apply the language-workbook edit: append new dialog rows (overworld_2 .. colony_4)
# to the 'en' worksheet, mirroring the Key/Value pairs added in the diff.
#
# Cells are written in the same order the original author typed them (column A
# for a block of rows, then column B for that same block) so the shared-string
# table ends up in the same sequence as the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column letter, text
$writes = @(
    130, 'A', 'overworld_2_intro_0'
    131, 'A', 'overworld_2_intro_1'
    131, 'B', 'In that case, we’ll need to check the wind readings of Earth.'
    132, 'A', 'overworld_2_wind_0'
    133, 'A', 'overworld_2_wind_temp_0'
    134, 'A', 'overworld_2_wind_temp_1'
    135, 'A', 'overworld_2_wind_temp_2'
    135, 'B', 'This cycle continues as the wind travels, building up more speed along the way.'
    136, 'A', 'overworld_2_post_intro_0'
    137, 'A', 'overworld_2_post_intro_1'
    138, 'A', 'overworld_2_post_intro_2'
    136, 'B', 'Anyhow, let’s find places where the wind might by strong. Try looking for areas prone to hurricanes along coastal regions.'
    137, 'B', 'Remember to check the different seasons to see how the wind changes to various positions across Earth.'
    138, 'B', 'This time around, there are more than one hotspot to discover on the map. Only one of them is will match with the frogs’ criteria.'
    130, 'B', 'Our next batch of frogs are keen on living in a windy environment, as well as hot and humid.'
    132, 'B', 'Notice the winds forming into a swirly motion? These are tropical cyclones, sometimes referred to as: hurricanes, or typhoons.'
    133, 'B', 'As you can see, the warm energy from the ocean rises, forming low pressure from below. This causes more air to fill in.'
    134, 'B', 'The air filling in must then rise because of the heat, accumulating water which release more heat. '
    139, 'A', 'colony_2_intro_0'
    140, 'A', 'colony_2_intro_1'
    141, 'A', 'colony_2_intro_2'
    139, 'B', 'Looks like we’ve landed in a tropical climate. Where it’s hot and humid all year round with plenty of rain.'
    140, 'B', 'Though it looks like we’re getting more rain than usual, something is afoot...'
    141, 'B', 'Perhaps we should take a careful look at the weather forecast.'
    142, 'A', 'colony_2_mushroom_0'
    143, 'A', 'colony_2_mushroom_1'
    144, 'A', 'colony_2_mushroom_2'
    142, 'B', 'Uh oh, a mushroom has grown near one of our structures!'
    143, 'B', 'Since there''s a lot of moisture in the region, the fungi that grow these mushrooms are able to absorb a lot of nutrients.'
    144, 'B', 'Their spores appear to be harmful to the frogs! Make sure to have a gardener around to take care of them!'
    145, 'A', 'colony_2_fly_0'
    146, 'A', 'colony_2_fly_1'
    147, 'A', 'colony_2_fly_2'
    145, 'B', 'Look out, it''s a beetle!'
    146, 'B', 'Due to the consistent warmth in tropical climates, insects are able to thrive throughout the year.'
    147, 'B', 'These troublesome insects can be dealt with by a hero frog. Make sure to have one around to rout them out.'
    148, 'A', 'colony_2_hazzard_0'
    149, 'A', 'colony_2_hazzard_1'
    150, 'A', 'colony_2_hazzard_2'
    151, 'A', 'colony_2_hazzard_3'
    152, 'A', 'colony_2_hazzard_4'
    148, 'B', 'Take cover, a hurricane is heading our way!'
    149, 'B', 'As mentioned before about hurricanes: the wind speed that has accumulated over low pressure from the surface has reached critical speed.'
    150, 'B', 'Our frogs must take cover. Fortunately, their structures are made of sturdy stuff, causing it to withstand the staggering winds!'
    151, 'B', 'However, along with strong winds, the water that is released from the storm will cause flood across the land.'
    152, 'B', 'Make certain that no important structures are within the flooding area, or they will get damaged.'
    153, 'A', 'overworld_3_intro_0'
    154, 'A', 'overworld_3_intro_1'
    153, 'B', 'Our next batch of frogs are looking for a warm place with low humidity, and some nice breeze.'
    154, 'B', 'In that case, we should look for a desert climate!'
    155, 'A', 'overworld_3_investigate_0'
    156, 'A', 'overworld_3_investigate_1'
    157, 'A', 'overworld_3_investigate_2'
    155, 'B', 'Now it may seem that these frogs would want to bask in the sun all day long. '
    156, 'B', 'However, they will still need some water source to sustain themselves!'
    157, 'B', 'Look for a suitable place where there are underground waters that we can extract.'
    158, 'A', 'colony_3_intro_0'
    159, 'A', 'colony_3_intro_1'
    158, 'B', 'The desert climate...This dry and hot environment will leave us with little to no water for our plants.'
    159, 'B', 'Fortunately, there are underground waters we can extract from to make this land more habitable.'
    160, 'A', 'colony_3_water_0'
    161, 'A', 'colony_3_water_1'
    160, 'B', 'Since the ground here is not ideal for growing plants, we will have to do a bit of landscaping.'
    161, 'B', 'First, we will need to build a water tank where water is accessible.'
    162, 'A', 'colony_3_landscaping_0'
    162, 'B', 'Now that we have a water source, summon a landscaper to irrigate the land.'
    163, 'A', 'colony_3_landscaping_complete_0'
    164, 'A', 'colony_3_landscaping_complete_1'
    163, 'B', 'Excellent! Now that there’s an irrigated land, you can now place a plant on it. Do this now.'
    164, 'B', 'We can proceed onward once we have increased the population.'
    165, 'A', 'overworld_4_intro_0'
    166, 'A', 'overworld_4_intro_1'
    165, 'B', 'This is our final batch of frogs, and they seem eager to settle in the highlands where it’s cold and snowy.'
    166, 'B', 'Why don’t we look for a spot in the mountainous area.'
    167, 'A', 'overworld_4_investigate_0'
    168, 'A', 'overworld_4_investigate_1'
    167, 'B', 'Although we are in an area that is mostly a tropical climate, remember that altitude can also affect the climate.'
    168, 'B', 'Go further up where the air pressure and temperature are lower.'
    169, 'A', 'colony_4_intro_0'
    170, 'A', 'colony_4_intro_1'
    169, 'B', 'The highland climate is quite comfy despite the consistent cold weather. We’ll need more than usual power to keep our houses warm.'
    170, 'B', 'Just like in the desert climate, the ground here is not ideal for plants to grow, so landscaping will be crucial.'
    171, 'A', 'colony_4_landscape_0'
    172, 'A', 'colony_4_landscape_1'
    173, 'A', 'colony_4_landscape_2'
    171, 'B', 'Since we can’t place plants on these rigid grounds, we’ll once again need the help of landscapers to shape the land.'
    172, 'B', 'You won’t have to worry about where to place the water thank this time around.'
    173, 'B', 'We will be able to proceed once the frog population has increased.'
    174, 'A', 'colony_4_cave_0'
    175, 'A', 'colony_4_cave_1'
    176, 'A', 'colony_4_cave_2'
    174, 'B', 'Uh oh! A cave has emerged from the ground!'
    175, 'B', 'Critters will keep emerging from these caves so long as it remains. Fortunately, an engineer can demolish it.'
    176, 'B', 'Make sure to also have a hero frog around to deal with the critters, while the engineers do their work!'
)

for ($i = 0; $i -lt $writes.Count; $i += 3) {
    $row = $writes[$i]
    $col = $writes[$i + 1]
    $text = $writes[$i + 2]
    if ($col -eq 'A') {
        $ws.Cells.Item($row, 1).Value = $text
    } else {
        $ws.Cells.Item($row, 2).Value = $text
    }
}

# Rows whose Value (column B) cell uses the vertical-center style, matching
# the existing sheet's alternating format.
$centeredRows = @(133, 135, 137, 138, 146, 152, 153, 158, 161, 162, 163, 164, 165, 166, 167, 168, 172, 173, 174, 175, 176)
foreach ($row in $centeredRows) {
    $ws.Cells.Item($row, 2).VerticalAlignment = -4108
}

# Move the saved selection to mirror the author's last cursor spot.
$ws.Activate()
$ws.Range("A174").Select()

Write-Host "Added $($centeredRows.Count) centered rows; wrote $($writes.Count / 3) cells to sheet $($ws.Name)."
